$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume figures (and the Arweave/Monero row swap).
# Each target cell is briefly marked as Text ("@") before the new value is written
# so that Excel does not auto-convert numeric-looking strings (e.g. "425.97") into
# real numbers; the temporary formatting is cleared immediately afterwards so the
# cell ends up with the same (default/unstyled) appearance it started with.
function Set-TextValue([string]$cellAddr, [string]$value) {
    $cell = $ws.Range($cellAddr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue "D2" "63.257.58"
Set-TextValue "E2" "  -0.10%  "
Set-TextValue "D3" "3.297.00"
Set-TextValue "E3" "  +1.46%  "
Set-TextValue "E4" "  -0.05%  "
Set-TextValue "D5" "601.38"
Set-TextValue "E5" "  +1.23%  "
Set-TextValue "D6" "139.61"
Set-TextValue "E6" "  -0.52%  "
Set-TextValue "D8" "3.296.06"
Set-TextValue "E8" "  +1.61%  "
Set-TextValue "E9" "  -0.58%  "
Set-TextValue "E10" "  +0.98%  "
Set-TextValue "D11" "5.47"
Set-TextValue "E11" "  +1.85%  "
Set-TextValue "D12" "0.465"
Set-TextValue "E12" "  +0.07%  "
Set-TextValue "D14" "34.35"
Set-TextValue "E14" "  +0.41%  "
Set-TextValue "D15" "3.839.60"
Set-TextValue "E15" "  +1.41%  "
Set-TextValue "E16" "  +1.26%  "
Set-TextValue "D17" "3.297.84"
Set-TextValue "E17" "  +1.53%  "
Set-TextValue "D18" "63.352.06"
Set-TextValue "E18" "  -0.04%  "
Set-TextValue "E19" "  +0.73%  "
Set-TextValue "D20" "475.20"
Set-TextValue "E20" "  +0.34%  "
Set-TextValue "D21" "13.90"
Set-TextValue "E21" "  -1.87%  "
Set-TextValue "E22" "  +0.12%  "
Set-TextValue "D23" "7.91"
Set-TextValue "E23" "  -0.22%  "
Set-TextValue "D24" "13.76"
Set-TextValue "E24" "  +4.66%  "
Set-TextValue "D25" "84.87"
Set-TextValue "E25" "  +1.33%  "
Set-TextValue "E26" "  +0.03%  "
Set-TextValue "E27" "  +0.83%  "
Set-TextValue "E28" "  -0.12%  "
Set-TextValue "D29" "7.11"
Set-TextValue "E29" "  -1.68%  "
Set-TextValue "D30" "8.07"
Set-TextValue "E30" "  +0.03%  "
Set-TextValue "D32" "28.41"
Set-TextValue "E32" "  +2.97%  "
Set-TextValue "E33" "  -2.39%  "
Set-TextValue "E34" "  -1.20%  "
Set-TextValue "E35" "  -0.17%  "
Set-TextValue "E36" "  +1.12%  "
Set-TextValue "D37" "52.18"
Set-TextValue "E37" "  -0.79%  "
Set-TextValue "E38" "  +2.57%  "
Set-TextValue "E39" "  +1.37%  "
Set-TextValue "D40" "3.110.89"
Set-TextValue "E40" "  +4.22%  "
Set-TextValue "D41" "425.97"
Set-TextValue "D42" "0.118"
Set-TextValue "E42" "  +7.22%  "
Set-TextValue "E43" "  -1.09%  "
Set-TextValue "E44" "  -1.93%  "
Set-TextValue "E45" "  -1.57%  "
Set-TextValue "E46" "  +1.18%  "
Set-TextValue "B47" "Monero"
Set-TextValue "C47" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D47" "128.15"
Set-TextValue "E47" "  +4.45%  "
Set-TextValue "B48" "Arweave"
Set-TextValue "C48" "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D48" "36.00"
Set-TextValue "E48" "  +8.02%  "
Set-TextValue "E50" "  +1.20%  "
Set-TextValue "E51" "  -1.37%  "
